# Update column G ("K") values on Sheet1 from the regenerated save_data
# (Strike# replaced by K). Only the rows whose K value actually changed
# are touched; all other cells are left exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 2
    3  = 0
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 1
    31 = 0
    32 = 0
    34 = 2
    35 = 1
    36 = 1
    37 = 0
    38 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
